# Codes update (apr 26 deployment)
# Mirrors a refreshed "DSDV-INFC-1601_qsrv1 NIP_INSIDENIP tblCVXCodes" query:
#   - table/defined-name renamed from "..._1" to the un-suffixed name
#   - table range grows from A1:H239 to A1:H247 (8 new CVX codes)
#   - one existing date (H235) is corrected
#   - the update_date column's custom date format drops the leading zero (dd -> d)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WEB cvx list")

function Set-TextCell($row, $col, $val) {
    # Force the cell to be stored as text even when the value looks numeric
    # ("512") or boolean-like ("False"), matching the source query column type.
    $rng = $ws.Cells.Item($row, $col)
    $rng.Value = "'" + $val
    $rng.Style = "Normal"
}

# ---------------------------------------------------------------------------
# 1) Fix an existing row's update_date (CVX 257 -> 16-Mar-2022)
# ---------------------------------------------------------------------------
$ws.Cells.Item(235, 8).Value = 44636

# ---------------------------------------------------------------------------
# 2) Append the 8 new CVX code rows (240-247) pulled in by the refresh
# ---------------------------------------------------------------------------
# Row 240
Set-TextCell 240 1 "512"
$ws.Cells.Item(240, 2).Value = "SARS-COV-2 COVID-19 VLP Non-US Vaccine (Medicago, Covifenz)"
$ws.Cells.Item(240, 3).Value = "SARS-COV-2 COVID-19 Virus Like Particle (VLP) Non-US Vaccine Product (Medicago, Covifenz)"
$ws.Cells.Item(240, 4).Value = "Pandemic Non-US Vaccine not Authorized by WHO - ACIP does recognize towards immunity in US"
$ws.Cells.Item(240, 5).Value = "Non-US"
$ws.Cells.Item(240, 6).Value = 262
Set-TextCell 240 7 "False"
$ws.Cells.Item(240, 8).Value = 44663

# Row 241
Set-TextCell 241 1 "513"
$ws.Cells.Item(241, 2).Value = "SARS-COV-2 COVID-19 PS Non-US Vaccine (Anhui Zhifei Longcom, Zifivax)"
$ws.Cells.Item(241, 3).Value = "SARS-COV-2 COVID-19 Protein Subunit Non-US Vaccine Product (Anhui Zhifei Longcom, Zifivax)"
$ws.Cells.Item(241, 4).Value = "Pandemic Non-US Vaccine not Authorized by WHO - not counted toward immunity in US"
$ws.Cells.Item(241, 5).Value = "Non-US"
$ws.Cells.Item(241, 6).Value = 263
Set-TextCell 241 7 "False"
$ws.Cells.Item(241, 8).Value = 44663

# Row 242
Set-TextCell 242 1 "514"
$ws.Cells.Item(242, 2).Value = "SARS-COV-2 COVID-19 DNA Non-US Vaccine (Zydus Cadila, ZyCoV-D)"
$ws.Cells.Item(242, 3).Value = "SARS-COV-2 COVID-19 DNA Non-US Vaccine Product (Zydus Cadila, ZyCoV-D)"
$ws.Cells.Item(242, 4).Value = "Pandemic Non-US Vaccine not Authorized by WHO - not counted toward immunity in US"
$ws.Cells.Item(242, 5).Value = "Non-US"
$ws.Cells.Item(242, 6).Value = 264
Set-TextCell 242 7 "False"
$ws.Cells.Item(242, 8).Value = 44663

# Row 243
Set-TextCell 243 1 "515"
$ws.Cells.Item(243, 2).Value = "SARS-COV-2 COVID-19 PS Non-US Vaccine (Medigen, MVC-COV1901)"
$ws.Cells.Item(243, 3).Value = "SARS-COV-2 COVID-19 Protein Subunit Non-US Vaccine Product (Medigen, MVC-COV1901)"
$ws.Cells.Item(243, 4).Value = "Pandemic Non-US Vaccine not Authorized by WHO - not counted toward immunity in US"
$ws.Cells.Item(243, 5).Value = "Non-US"
$ws.Cells.Item(243, 6).Value = 265
Set-TextCell 243 7 "False"
$ws.Cells.Item(243, 8).Value = 44663

# Row 244
Set-TextCell 244 1 "516"
$ws.Cells.Item(244, 2).Value = "COV-2 COVID-19 Inactivated Non-US Vaccine Product (Minhai Biotechnology Co, KCONVAC)"
$ws.Cells.Item(244, 3).Value = "SARS-COV-2 COVID-19 Inactivated Non-US Vaccine Product (Minhai Biotechnology Co, KCONVAC)"
$ws.Cells.Item(244, 4).Value = "Pandemic Non-US Vaccine not Authorized by WHO - not counted toward immunity in US"
$ws.Cells.Item(244, 5).Value = "Non-US"
$ws.Cells.Item(244, 6).Value = 266
Set-TextCell 244 7 "False"
$ws.Cells.Item(244, 8).Value = 44663

# Row 245
Set-TextCell 245 1 "517"
$ws.Cells.Item(245, 2).Value = "SARS-COV-2 COVID-19 PS Non-US Vaccine (Biological E Limited, Corbevax)"
$ws.Cells.Item(245, 3).Value = "SARS-COV-2 COVID-19 Protein Subunit Non-US Vaccine Product (Biological E Limited, Corbevax)"
$ws.Cells.Item(245, 4).Value = "Pandemic Non-US Vaccine not Authorized by WHO - not counted toward immunity in US"
$ws.Cells.Item(245, 5).Value = "Non-US"
$ws.Cells.Item(245, 6).Value = 267
Set-TextCell 245 7 "False"
$ws.Cells.Item(245, 8).Value = 44663

# Row 246
Set-TextCell 246 1 "227"
$ws.Cells.Item(246, 2).Value = "COVID-19, mRNA, LNP-S, PF, pediatric 50 mcg/0.5 mL dose"
$ws.Cells.Item(246, 3).Value = "SARS-COV-2 (COVID-19) vaccine, mRNA, spike protein, LNP, preservative free, pediatric 50 mcg/0.5 mL dose"
$ws.Cells.Item(246, 4).Value = "Pre-EUA Moderna Pediatric 6yr to<12 yr vaccine 2.5 mL vial, 50 mcg/0.5 mL dose"
$ws.Cells.Item(246, 5).Value = "Active"
$ws.Cells.Item(246, 6).Value = 268
Set-TextCell 246 7 "False"
$ws.Cells.Item(246, 8).Value = 44663

# Row 247
Set-TextCell 247 1 "228"
$ws.Cells.Item(247, 2).Value = "COVID-19, mRNA, LNP-S, PF, pediatric 25 mcg/0.25 mL dose"
$ws.Cells.Item(247, 3).Value = "SARS-COV-2 (COVID-19) vaccine, mRNA, spike protein, LNP, preservative free, pediatric 25 mcg/0.25 mL dose"
$ws.Cells.Item(247, 4).Value = "Pre-EUA Moderna Pediatric 6mo to<6yr 2.5 mL vial, 25 mcg/0.25 mL dose"
$ws.Cells.Item(247, 5).Value = "Active"
$ws.Cells.Item(247, 6).Value = 269
Set-TextCell 247 7 "False"
$ws.Cells.Item(247, 8).Value = 44663

# ---------------------------------------------------------------------------
# 3) Grow + rename the query table and its backing defined name
#    (Excel drops the "_1" disambiguation suffix once it is unique again)
# ---------------------------------------------------------------------------
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:H247"))
$tbl.Name = "Table_DSDV_INFC_1601_qsrv1_NIP_INSIDENIP_tblCVXCodes"

$defName = $wb.Names.Item(1)
$defName.Name = "DSDV_INFC_1601_qsrv1_NIP_INSIDENIP_tblCVXCodes"
$defName.RefersTo = "='WEB cvx list'!`$A`$1:`$H`$247"

# ---------------------------------------------------------------------------
# 4) Drop the leading zero from the update_date custom format
#    ("01-Jan-22" -> "1-Jan-22") across the whole column
# ---------------------------------------------------------------------------
$ws.Range("H2:H247").NumberFormat = "[$-409]d\-mmm\-yy;@"

# ---------------------------------------------------------------------------
# 5) Cosmetic refresh artifacts: column widths, page scale, view/selection
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 30.11
$ws.Columns.Item(5).ColumnWidth = 15.89
$ws.Columns.Item(8).ColumnWidth = 12.78

$ws.PageSetup.Zoom = 19

$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B4").Select()
